$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'243.79"
$ws.Range("D4").Formula = "'5.395"
$ws.Range("D5").Formula = "'0.05936"
$ws.Range("D6").Formula = "'3.435"
$ws.Range("D7").Formula = "'6.507"
$ws.Range("D8").Formula = "'0.8114"
$ws.Range("D9").Formula = "'0.9284"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Formula = "'0.1434"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Formula = "'0.07422"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Formula = "'0.03234"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Formula = "'0.03076"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Formula = "'0.09357"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Formula = "'3.851"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Formula = "'0.001580"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Formula = "'0.04700"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Formula = "'0.0005987"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("D19").Formula = "'0.005942"
$ws.Range("D20").Formula = "'0.001257"
$ws.Range("E20").Value = "19BitKanKANBestin24h"
$ws.Range("D21").Formula = "'0.004787"
$ws.Range("D23").Formula = "'3.560"
$ws.Range("D26").Formula = "'0.1332"
$ws.Range("D27").Formula = "'0.0002342"
$ws.Range("E27").Value = "26UpBotsUBXT"
$ws.Range("D40").Formula = "'0.03922"
$ws.Range("D41").Formula = "'0.003073"
$ws.Range("E41").Value = "40KickTokenKICKWorstin24h"
$ws.Range("D42").Formula = "'0.1075"
$ws.Range("D43").Formula = "'0.002590"
$ws.Range("D44").Formula = "'0.008215"
$ws.Range("D45").Formula = "'0.00005187"
$ws.Range("D46").Formula = "'0.00000000751"
$ws.Range("D47").Formula = "'0.6668"
$ws.Range("D48").Formula = "'0.002067"
$ws.Range("D49").Formula = "'0.00002103"
$ws.Range("D50").Formula = "'0.0002002"
